$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.182.09"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "3.073.55"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'522.31"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").Value = "'135.53"
$ws.Range("E6").Value = "  -4.89%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.074.00"
$ws.Range("E8").Value = "  -1.69%  "
$ws.Range("E9").Value = "  +4.99%  "
$ws.Range("D10").Value = "'7.30"
$ws.Range("E10").Value = "  +1.88%  "
$ws.Range("E11").Value = "  -2.76%  "
$ws.Range("D12").Value = "'0.401"
$ws.Range("E12").Value = "  +2.16%  "
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("D14").Value = "3.595.69"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").Value = "'25.19"
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").Value = "'0.0000161"
$ws.Range("E16").Value = "  -3.07%  "
$ws.Range("D17").Value = "57.234.74"
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").Value = "3.076.25"
$ws.Range("E18").Value = "  -1.94%  "
$ws.Range("D19").Value = "'5.87"
$ws.Range("E19").Value = "  -4.11%  "
$ws.Range("D20").Value = "'12.43"
$ws.Range("E20").Value = "  -3.12%  "
$ws.Range("E21").Value = "  -2.07%  "
$ws.Range("D22").Value = "'349.02"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "'68.97"
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("E25").Value = "  -3.06%  "
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").Value = "0.0₃0860"
$ws.Range("E28").Value = "  -7.56%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -2.69%  "
$ws.Range("D31").Value = "'1.86"
$ws.Range("E31").Value = "  -1.12%  "
$ws.Range("D32").Value = "'20.90"
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("D33").Value = "'5.77"
$ws.Range("E33").Value = "  -9.81%  "
$ws.Range("D34").Value = "'159.27"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").Value = "'4.82"
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("E36").Value = "  -4.81%  "
$ws.Range("E37").Value = "  -3.47%  "
$ws.Range("D38").Value = "'25.32"
$ws.Range("E38").Value = "  -3.89%  "
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("D40").Value = "'0.0654"
$ws.Range("E40").Value = "  -1.87%  "
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("E42").Value = "  -6.56%  "
$ws.Range("D43").Value = "'0.693"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D44").Value = "2.412.55"
$ws.Range("E44").Value = "  +6.54%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "3.112.56"
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("E49").Value = "  -2.61%  "
$ws.Range("D50").Value = "'0.935"
$ws.Range("E50").Value = "  -6.88%  "
$ws.Range("D51").Value = "'19.51"
$ws.Range("E51").Value = "  -5.56%  "
